$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1712.375
$ws.Range("I40").Value = 1588.7222
$ws.Range("K40").Value = 1588.7222
$ws.Range("M40").Value = -1413.7222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7689.3687
$ws.Range("I51").Value = 11809.1
$ws.Range("J51").Value = 3111.889
$ws.Range("K51").Value = 11809.1
$ws.Range("L51").Value = 3111.889
$ws.Range("M51").Value = -11325.1
$ws.Range("N51").Value = -4079.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8632
$ws.Range("I69").Value = 7000
$ws.Range("J69").Value = 9720
$ws.Range("K69").Value = 21000
$ws.Range("L69").Value = 29160
$ws.Range("M69").Value = -20126
$ws.Range("N69").Value = -30908

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 8632
$ws.Range("I72").Value = 7000
$ws.Range("J72").Value = 9720
$ws.Range("K72").Value = 63000
$ws.Range("L72").Value = 87480
$ws.Range("M72").Value = -58632
$ws.Range("N72").Value = -96216

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -24992

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 795.1667
$ws.Range("I99").Value = 447.5
$ws.Range("J99").Value = 1490.5
$ws.Range("K99").Value = 1342.5
$ws.Range("L99").Value = 4471.5
$ws.Range("M99").Value = 155.5
$ws.Range("N99").Value = -7467.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 429.33334
$ws.Range("I101").Value = 408
$ws.Range("J101").Value = 600
$ws.Range("K101").Value = 1224
$ws.Range("L101").Value = 1800
$ws.Range("M101").Value = 398
$ws.Range("N101").Value = -5044

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1179.5695
$ws.Range("I137").Value = 1029.0344
$ws.Range("J137").Value = 1281.093
$ws.Range("K137").Value = 3087.1032
$ws.Range("L137").Value = 3843.279
$ws.Range("M137").Value = -537.1032
$ws.Range("N137").Value = -8943.279

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4824.9453
$ws.Range("I138").Value = 3183
$ws.Range("J138").Value = 5333.1665
$ws.Range("K138").Value = 9549
$ws.Range("L138").Value = 15999.4995
$ws.Range("M138").Value = -4409
$ws.Range("N138").Value = -26279.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 16
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 13000
$ws.Range("J9").Value = 13000
$ws.Range("L9").Value = 13000
$ws.Range("N9").Value = -13340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 13000
$ws.Range("J20").Value = 13000
$ws.Range("L20").Value = 13000
$ws.Range("N20").Value = -13540

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2076.77
$ws.Range("I32").Value = 1765.0581
$ws.Range("J32").Value = 3991.5715
$ws.Range("K32").Value = 1765.0581
$ws.Range("L32").Value = 3991.5715
$ws.Range("M32").Value = -1478.0581
$ws.Range("N32").Value = -4565.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2128.6365
$ws.Range("I45").Value = 1827.1364
$ws.Range("K45").Value = 1827.1364
$ws.Range("M45").Value = -1450.1364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2249.7812
$ws.Range("I61").Value = 1503.8334
$ws.Range("J61").Value = 2421.923
$ws.Range("K61").Value = 1503.8334
$ws.Range("L61").Value = 2421.923
$ws.Range("M61").Value = -1291.8334
$ws.Range("N61").Value = -2845.923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 169196.67
$ws.Range("I102").Value = 335326.66
$ws.Range("K102").Value = 335326.66
$ws.Range("M102").Value = -333704.66

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2670.4482
$ws.Range("I122").Value = 3200.2778
$ws.Range("J122").Value = 1803.4546
$ws.Range("K122").Value = 9600.8334
$ws.Range("L122").Value = 5410.3638
$ws.Range("M122").Value = -7150.8334
$ws.Range("N122").Value = -10310.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2828.6482
$ws.Range("I132").Value = 3575
$ws.Range("J132").Value = 1743.0454
$ws.Range("K132").Value = 10725
$ws.Range("L132").Value = 5229.1362
$ws.Range("M132").Value = -8195
$ws.Range("N132").Value = -10289.1362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2249.7812
$ws.Range("I136").Value = 1503.8334
$ws.Range("J136").Value = 2421.923
$ws.Range("K136").Value = 4511.5002
$ws.Range("L136").Value = 7265.768999999999
$ws.Range("M136").Value = -1961.5002
$ws.Range("N136").Value = -12365.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21718
$ws.Range("I82").Value = 4131.2
$ws.Range("J82").Value = 36373.668
$ws.Range("K82").Value = 4131.2
$ws.Range("L82").Value = 36373.668
$ws.Range("M82").Value = -3748.2
$ws.Range("N82").Value = -37139.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 21718
$ws.Range("I85").Value = 4131.2
$ws.Range("J85").Value = 36373.668
$ws.Range("K85").Value = 4131.2
$ws.Range("L85").Value = 36373.668
$ws.Range("M85").Value = -2805.2
$ws.Range("N85").Value = -39025.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 6428
$ws.Range("I97").Value = 6428
$ws.Range("K97").Value = 6428
$ws.Range("M97").Value = -5437

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15578.411
$ws.Range("I31").Value = 33367.484
$ws.Range("J31").Value = 2448.3809
$ws.Range("K31").Value = 33367.484
$ws.Range("L31").Value = 2448.3809
$ws.Range("M31").Value = -33072.484
$ws.Range("N31").Value = -3038.3809

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 15578.411
$ws.Range("I34").Value = 33367.484
$ws.Range("J34").Value = 2448.3809
$ws.Range("K34").Value = 33367.484
$ws.Range("L34").Value = 2448.3809
$ws.Range("M34").Value = -33165.484
$ws.Range("N34").Value = -2852.3809

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 25000
$ws.Range("I47").Value = 15000
$ws.Range("J47").Value = 35000
$ws.Range("K47").Value = 15000
$ws.Range("L47").Value = 35000
$ws.Range("M47").Value = -14434
$ws.Range("N47").Value = -36132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 28000
$ws.Range("J48").Value = 28000
$ws.Range("L48").Value = 28000
$ws.Range("N48").Value = -28952

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 43000
$ws.Range("J112").Value = 43000
$ws.Range("L112").Value = 43000
$ws.Range("N112").Value = -45954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 69000
$ws.Range("J131").Value = 69000
$ws.Range("L131").Value = 69000
$ws.Range("N131").Value = -79080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3133.5
$ws.Range("I132").Value = 3003.476
$ws.Range("J132").Value = 3679.6
$ws.Range("K132").Value = 9010.428
$ws.Range("L132").Value = 11038.8
$ws.Range("M132").Value = -6480.428
$ws.Range("N132").Value = -16098.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1891.6666
$ws.Range("I134").Value = 1418.5454
$ws.Range("K134").Value = 4255.6362
$ws.Range("M134").Value = -1720.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 295
$ws.Range("I18").Value = 253.25
$ws.Range("J18").Value = 517.6667
$ws.Range("K18").Value = 759.75
$ws.Range("L18").Value = 1553.0001
$ws.Range("M18").Value = -590.75
$ws.Range("N18").Value = -1891.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1500.1666
$ws.Range("I114").Value = 1499.5
$ws.Range("J114").Value = 1500.5
$ws.Range("K114").Value = 4498.5
$ws.Range("L114").Value = 4501.5
$ws.Range("M114").Value = -1244.5
$ws.Range("N114").Value = -11009.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 3107.5
$ws.Range("I126").Value = 3476.6667
$ws.Range("K126").Value = 10430.0001
$ws.Range("M126").Value = -5490.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1235936.9
$ws.Range("J131").Value = 1450697.9
$ws.Range("L131").Value = 4352093.699999999
$ws.Range("N131").Value = -4362173.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3213.75
$ws.Range("I132").Value = 4333.3335
$ws.Range("J132").Value = 3016.1765
$ws.Range("K132").Value = 39000.0015
$ws.Range("L132").Value = 27145.5885
$ws.Range("M132").Value = -36470.0015
$ws.Range("N132").Value = -32205.5885

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5456152.5
$ws.Range("I12").Value = 5456152.5
$ws.Range("K12").Value = 5456152.5
$ws.Range("M12").Value = -5456012.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 47620572
$ws.Range("I97").Value = 66668210
$ws.Range("J97").Value = 1486.8334
$ws.Range("K97").Value = 66668210
$ws.Range("L97").Value = 1486.8334
$ws.Range("M97").Value = -66667714
$ws.Range("N97").Value = -2478.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3529.4443
$ws.Range("I102").Value = 4095.2856
$ws.Range("J102").Value = 1549
$ws.Range("K102").Value = 4095.2856
$ws.Range("L102").Value = 1549
$ws.Range("M102").Value = -2473.2856
$ws.Range("N102").Value = -4793

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4499.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3589.1428
$ws.Range("I132").Value = 3550
$ws.Range("K132").Value = 10650
$ws.Range("M132").Value = -8120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2535.2354
$ws.Range("I136").Value = 2078.1428
$ws.Range("K136").Value = 6234.428400000001
$ws.Range("M136").Value = -3684.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 9400
$ws.Range("I14").Value = 750
$ws.Range("J14").Value = 20933.334
$ws.Range("K14").Value = 750
$ws.Range("L14").Value = 20933.334
$ws.Range("M14").Value = -582
$ws.Range("N14").Value = -21269.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1176.5264
$ws.Range("I122").Value = 1284.75
$ws.Range("K122").Value = 3854.25
$ws.Range("M122").Value = -1404.25
